$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data (header + rows).
$lastRow = $ws.UsedRange.Rows.Count

# Swap the "category-code" (F) and "category-name" (G) columns, for the
# header row and every data row, including their header cells. Using
# Range.Copy() (rather than re-typing .Value) preserves each cell's
# original storage type, so numeric-looking codes such as "111" remain
# text instead of being reinterpreted as numbers.
$srcF = $ws.Range("F1:F$lastRow")
$srcG = $ws.Range("G1:G$lastRow")
$tmp  = $ws.Range("I1:I$lastRow")

$srcF.Copy($tmp)
$srcG.Copy($srcF)
$tmp.Copy($srcG)

$tmp.ClearContents()
